$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear header cells that no longer hold food labels (E1, G1)
$ws.Range("E1").Clear()
$ws.Range("G1").Clear()

# Set new / changed cell values
$ws.Range("C1").Value = "portion"
$ws.Range("N1").Value = "food"
$ws.Range("O1").Value = "food"
$ws.Range("P1").Value = "food"
$ws.Range("C2").Value = 200
$ws.Range("D2").Value = "banán"
$ws.Range("E2").Value = 300
$ws.Range("F2").Value = "zabkása"
$ws.Range("G2").Value = 200
$ws.Range("H2").Value = "tej"
$ws.Range("C3").Value = 300
$ws.Range("D3").Value = "csirke mell"
$ws.Range("E3").Value = 30
$ws.Range("F3").Value = "rizs"
$ws.Range("G3").Value = 200
$ws.Range("H3").Value = "uborka"
$ws.Range("C4").Value = 200
$ws.Range("D4").Value = "tojás"
$ws.Range("E4").Value = 300
$ws.Range("F4").Value = "kenyér"
$ws.Range("G4").Value = 100
$ws.Range("H4").Value = "paradicsom"
$ws.Range("C5").Value = 300
$ws.Range("D5").Value = "alma"
$ws.Range("E5").Value = 200
$ws.Range("F5").Value = "joghurt"
$ws.Range("C6").Value = 400
$ws.Range("D6").Value = "ponty"
$ws.Range("E6").Value = 400
$ws.Range("F6").Value = "burgonya"
$ws.Range("G6").Value = 200
$ws.Range("H6").Value = "saláta"
$ws.Range("C7").Value = 200
$ws.Range("D7").Value = "sajt"
$ws.Range("E7").Value = 200
$ws.Range("F7").Value = "uborka"
$ws.Range("C8").Value = 400
$ws.Range("D8").Value = "alma"
$ws.Range("E8").Value = 200
$ws.Range("F8").Value = "joghurt"
$ws.Range("C9").Value = 300
$ws.Range("D9").Value = "sonka"
$ws.Range("E9").Value = 300
$ws.Range("F9").Value = "burgonya"
$ws.Range("G9").Value = 300
$ws.Range("H9").Value = "saláta"
$ws.Range("C10").Value = 300
$ws.Range("D10").Value = "sonka"
$ws.Range("E10").Value = 300
$ws.Range("F10").Value = "uborka"
$ws.Range("C11").Value = 400
$ws.Range("D11").Value = "alma"
$ws.Range("C12").Value = 400
$ws.Range("D12").Value = "alma"
$ws.Range("C13").Value = 300
$ws.Range("D13").Value = "paradicsom"

# Update selection to match the saved view state
$ws.Range("C16").Select()
